$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Sheet "New Horizons": append rows 67-69
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("New Horizons")

$rows1 = @(
    @{ A="MEEPLE"; B="SQUEAK"; C="BONNIE"; D="JAE-YONG"; E="CHARLIE"; F="DARRYL"; G="Equipo 2"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T212546.000Z" },
    @{ A="MEEPLE"; B="SQUEAK"; C="BONNIE"; D="JAE-YONG"; E="CHARLIE"; F="DARRYL"; G="Equipo 2"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T212413.000Z" },
    @{ A="MEEPLE"; B="SQUEAK"; C="BONNIE"; D="JAE-YONG"; E="CHARLIE"; F="DARRYL"; G="Equipo 1"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T212118.000Z" }
)

$lastRow1 = 66
$equipo1Row1 = 66   # existing row formatted as "Equipo 1"
$equipo2Row1 = 60   # existing row formatted as "Equipo 2"

$startRow1 = 67
for ($i = 0; $i -lt $rows1.Count; $i++) {
    $r = $startRow1 + $i
    $row = $rows1[$i]

    # Copy full-row formatting from the last existing data row.
    $ws1.Range("A$lastRow1`:N$lastRow1").Copy()
    $ws1.Range("A$r`:N$r").PasteSpecial($xlPasteFormats)

    # Column G's fill depends on which team won; reuse the right sample cell.
    if ($row.G -eq "Equipo 1") {
        $ws1.Range("G$equipo1Row1").Copy()
    } else {
        $ws1.Range("G$equipo2Row1").Copy()
    }
    $ws1.Range("G$r").PasteSpecial($xlPasteFormats)

    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C
    $ws1.Cells.Item($r, 4).Value = $row.D
    $ws1.Cells.Item($r, 5).Value = $row.E
    $ws1.Cells.Item($r, 6).Value = $row.F
    $ws1.Cells.Item($r, 7).Value = $row.G
    $ws1.Cells.Item($r, 8).Value = $row.H
    $ws1.Cells.Item($r, 9).Value = $row.I
    $ws1.Cells.Item($r, 10).Value = $row.J
    $ws1.Cells.Item($r, 11).Value = $row.K
    $ws1.Cells.Item($r, 12).Value = $row.L
    $ws1.Cells.Item($r, 13).Value = $row.M
    $ws1.Cells.Item($r, 14).Value = $row.N
}

# ---------------------------------------------------------------------------
# Sheet "Layer Cake": append rows 73-76
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Layer Cake")

$rows2 = @(
    @{ A="KIT"; B="MEEPLE"; C="CHARLIE"; D="JUJU"; E="CORDELIUS"; F="GENE"; G="Equipo 1"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T211333.000Z" },
    @{ A="KIT"; B="MEEPLE"; C="CHARLIE"; D="JUJU"; E="CORDELIUS"; F="GENE"; G="Equipo 1"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T211114.000Z" },
    @{ A="KIT"; B="MEEPLE"; C="CHARLIE"; D="JUJU"; E="CORDELIUS"; F="GENE"; G="Equipo 2"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T210857.000Z" },
    @{ A="BERRY"; B="SQUEAK"; C="CROW"; D="BEA"; E="MR. P"; F="OLLIE"; G="Equipo 2"; H="NHG|Xemp"; I="KCP|Fade"; J="KCP|Tyrant"; K="TRB|Zeus 解開"; L="TRB|Lxffy"; M="TRB|R B M"; N="20250724T210254.000Z" }
)

$lastRow2 = 72
$equipo1Row2 = 64   # existing row formatted as "Equipo 1"
$equipo2Row2 = 72   # existing row formatted as "Equipo 2"

$startRow2 = 73
for ($i = 0; $i -lt $rows2.Count; $i++) {
    $r = $startRow2 + $i
    $row = $rows2[$i]

    $ws2.Range("A$lastRow2`:N$lastRow2").Copy()
    $ws2.Range("A$r`:N$r").PasteSpecial($xlPasteFormats)

    if ($row.G -eq "Equipo 1") {
        $ws2.Range("G$equipo1Row2").Copy()
    } else {
        $ws2.Range("G$equipo2Row2").Copy()
    }
    $ws2.Range("G$r").PasteSpecial($xlPasteFormats)

    $ws2.Cells.Item($r, 1).Value = $row.A
    $ws2.Cells.Item($r, 2).Value = $row.B
    $ws2.Cells.Item($r, 3).Value = $row.C
    $ws2.Cells.Item($r, 4).Value = $row.D
    $ws2.Cells.Item($r, 5).Value = $row.E
    $ws2.Cells.Item($r, 6).Value = $row.F
    $ws2.Cells.Item($r, 7).Value = $row.G
    $ws2.Cells.Item($r, 8).Value = $row.H
    $ws2.Cells.Item($r, 9).Value = $row.I
    $ws2.Cells.Item($r, 10).Value = $row.J
    $ws2.Cells.Item($r, 11).Value = $row.K
    $ws2.Cells.Item($r, 12).Value = $row.L
    $ws2.Cells.Item($r, 13).Value = $row.M
    $ws2.Cells.Item($r, 14).Value = $row.N
}
